$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E5").Value = "S"
$ws.Range("F5").Value = "logo en info over de museum hebben"
$ws.Range("F5").Select()
